# Updated Results with corrected code
# Applies corrected values to Sheet1 of 2050_ES.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): corrected Iron & steel demand, clear the
# Non-metallic minerals figure (becomes blank).
$ws.Range("B3").Value = 8751638.252402626
$ws.Range("D3").Value = ""

# Row 4 (Methanol): corrected Chemicals figure.
$ws.Range("C4").Value = 66.44873649342436

# Row 5 (Ammonia): corrected Chemicals figure.
$ws.Range("C5").Value = 1872.005399194904

# Row 7: relabel as "Biogas" and update its corrected
# Non-metallic minerals figure.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 5165.484642980836

# New row 8: re-introduce an "Other" row below Biogas, copying the
# label formatting used by the other row headers in column A.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 7014.255481548714
